# Slide 3 ("DDS Project1 - Reed_Miller.pptx", sldId 261): the text box
# "TextBox 7" has a bullet that currently ends with a trailing period:
#   "May need to be filtered out later to avoid misleading results."
# Remove the trailing period so it reads:
#   "May need to be filtered out later to avoid misleading results"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(2)

$textRange = $shp.TextFrame.TextRange
$targetParagraph = $textRange.Paragraphs(7, 1)
$targetParagraph.Text = "May need to be filtered out later to avoid misleading results"
